$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.845.17'
$ws.Range("E2").Value = '  -1.06%  '
$ws.Range("D3").Value = '3.424.12'
$ws.Range("E3").Value = '  +4.06%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '256.37'
$ws.Range("E5").Value = '  +1.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '656.45'
$ws.Range("E6").Value = '  +5.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.48'
$ws.Range("E7").Value = '  +2.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.429'
$ws.Range("E8").Value = '  +5.29%  '
$ws.Range("E9").Value = '  +8.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.998'
$ws.Range("E10").Value = '  -0.09%  '
$ws.Range("D11").Value = '3.420.95'
$ws.Range("E11").Value = '  +4.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.213'
$ws.Range("E12").Value = '  +6.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.10'
$ws.Range("E13").Value = '  +6.69%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.38'
$ws.Range("E14").Value = '  +16.40%  '
$ws.Range("E15").Value = '  +4.11%  '
$ws.Range("D16").Value = '97.671.30'
$ws.Range("E16").Value = '  -0.93%  '
$ws.Range("D17").Value = '4.055.66'
$ws.Range("E17").Value = '  +3.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.57'
$ws.Range("E18").Value = '  +35.01%  '
$ws.Range("D19").Value = '3.415.22'
$ws.Range("E19").Value = '  +3.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.53'
$ws.Range("E20").Value = '  +13.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.503'
$ws.Range("E21").Value = '  +53.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.08'
$ws.Range("E22").Value = '  +18.57%  '
$ws.Range("E23").Value = '  -0.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '507.94'
$ws.Range("E24").Value = '  +4.25%  '
$ws.Range("E25").Value = '  +2.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.09'
$ws.Range("E26").Value = '  +7.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '98.38'
$ws.Range("E27").Value = '  +10.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.87'
$ws.Range("E28").Value = '  +6.92%  '
$ws.Range("D29").Value = '3.600.45'
$ws.Range("E29").Value = '  +3.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.154'
$ws.Range("E30").Value = '  +10.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.43'
$ws.Range("E31").Value = '  +10.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.197'
$ws.Range("E32").Value = '  +3.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.26%  '
$ws.Range("B35").Value = 'PolygonEcosystemToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.573'
$ws.Range("E35").Value = '  +19.89%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '29.98'
$ws.Range("E36").Value = '  +7.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.24'
$ws.Range("E37").Value = '  +15.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.81'
$ws.Range("E38").Value = '  +8.20%  '
$ws.Range("E39").Value = '  +4.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.41'
$ws.Range("E40").Value = '  +14.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '516.43'
$ws.Range("E41").Value = '  +5.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '24.73'
$ws.Range("E42").Value = '  -0.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.862'
$ws.Range("E43").Value = '  +10.41%  '
$ws.Range("E44").Value = '  +2.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0414'
$ws.Range("E45").Value = '  +22.98%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.30'
$ws.Range("E46").Value = '  +6.03%  '
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.48'
$ws.Range("E47").Value = '  +15.60%  '
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.19'
$ws.Range("E49").Value = '  +11.92%  '
$ws.Range("E50").Value = '  +16.17%  '
$ws.Range("E51").Value = '  +6.69%  '
